$wb = $excel.ActiveWorkbook

# ============================================================
# 1) Add a new "2022-Q1" sheet (fund snapshot), positioned right
#    before the "总计" (totals) sheet. We duplicate "2021-Q4" so
#    the new sheet starts out with identical layout/styling, then
#    rename it and overwrite its data with the 2022-Q1 numbers.
# ============================================================
$srcSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$srcSheet.Copy($totalSheet)

# Re-resolve "总计" by name (its Index shifted because the copy was
# inserted right before it) and grab the sheet that landed just ahead
# of it - that is the duplicate we just created.
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Item($totalSheet.Index - 1)
$newSheet.Name = "2022-Q1"

# Fund-code / name / size / position columns (B:G) need to stay
# TEXT (they hold values like "008763" and "92.10" where leading
# and trailing zeros matter), so force a text number format before
# writing them - otherwise Excel auto-coerces them to numbers.
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("B2").Value = "008763"
$newSheet.Range("C2").Value = "天弘越南市场股票（QDII）A"
$newSheet.Range("D2").Value = "37.53"
$newSheet.Range("E2").Value = "92.10"
$newSheet.Range("F2").Value = "4.62"
$newSheet.Range("G2").Value = "1.7339"
$newSheet.Range("H2").Value = 7

$newSheet.Range("B3").Value = "008764"
$newSheet.Range("C3").Value = "天弘越南市场股票（QDII）C"
$newSheet.Range("D3").Value = "14.26"
$newSheet.Range("E3").Value = "92.10"
$newSheet.Range("F3").Value = "4.62"
$newSheet.Range("G3").Value = "0.6588"
$newSheet.Range("H3").Value = 7

# ============================================================
# 2) Update the "总计" (totals) sheet: push its three existing
#    quarterly summary rows down by one row and add a brand new
#    top row for 2022-Q1. We shift data with Copy + PasteSpecial
#    (values, then formats) rather than plain Range assignment so
#    each destination cell's style (e.g. the centered index-column
#    style on column A) is carried along correctly.
# ============================================================
$ws = $wb.Worksheets.Item("总计")

# row4 -> row5
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4163)
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)

# row3 -> row4
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4163)
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)

# row2 -> row3
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial(-4163)
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)

# new top data row: 2022-Q1
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q1"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 2.39

# re-number the shifted rows' index column
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
